$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "1.010", "0.00000000108") stay as text
# instead of being auto-converted to numbers by Excel.
$textCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D16","D17","D19","D21","D22","D24","D25","D26","D27","D28","D29","D31","D32","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.419.56"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.696.98"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "218.43"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "0.5480"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.2740"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "0.06447"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "0.07673"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "1.693.94"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "4.557"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "0.5849"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "65.65"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "26.467.50"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "191.56"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "6.268"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "148.91"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "0.1313"
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("D26").Value = "7.922"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("D27").Value = "15.84"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "0.06233"
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "3.610"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "3.596"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "0.6156"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "2.412"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "2.760"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").Value = "0.01656"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").Value = "1.118.66"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "6.105"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("D41").Value = "0.8809"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "101.20"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "1.847.06"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "57.66"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000108"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "8.196"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "0.05285"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "6.125"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "0.4302"
$ws.Range("E51").Value = "  +0.05%  "
